# 10.3.1 workbook update:
# The two section-header rows ("Age" / "Education") are reworded so the
# Kyrgyz/Russian/English captions read "By age (in years)" / "By education"
# style phrasing instead of the old bare "Age (in years)" / "Education"
# wording. Updating these six cells (A19:C19 and A29:C29) causes the old,
# now-unused shared-string entries to be dropped and new ones appended when
# the workbook is saved, mirroring the upstream XML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Kyrgyz captions (column A): age-group header, then education header
$ws.Range("A19").Value = "Жаш курагы боюнча (жылдарда)"
$ws.Range("A29").Value = "Билими боюнча"

# Russian captions (column B): age-group header, then education header
$ws.Range("B19").Value = "По возрасту (в годах)"
$ws.Range("B29").Value = "По образованию"

# English captions (column C): age-group header, then education header
$ws.Range("C19").Value = "By age (in years) "
$ws.Range("C29").Value = "By education"
